$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Fill in the Name property value (row 4, "Name" in A4) which was previously empty
$ws.Range("B4").Value = "DescnonqualifiantVs"

# Update the Date property value (row 8, "Date" in A8)
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
